$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("H1").Value = "Date"

# Date values (as serials, so no time-of-day component sneaks in)
$ws.Range("H2").Value = 42370
$ws.Range("H3").Value = 43092
$ws.Range("H4").Value = 25934
$ws.Range("H5").Value = 77514

# Apply the same custom date format to the whole new column
$ws.Range("H1:H5").NumberFormat = "D/\ MMMM\ YYYY"

# Reset the selection back to the top-left cell
[void]$ws.Range("A1").Select()
